$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "yearly_income < 50,000 AND" -> "yearly_income < 50,000 "
#   (drop the trailing "AND", keep the trailing space)
# -----------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute(
    " < 50,000 AND", $true, $false, $false, $false, $false,
    $true, 1, $false, " < 50,000 ", 2)

# -----------------------------------------------------------------
# Change 2: split the "::Ask GPA only if the student is new::" run
# so the leading "::" is wrapped in gramStart/gramEnd proofErr marks,
# and the preceding tab becomes its own run (no longer sharing a run
# with the text).
# -----------------------------------------------------------------
$para = $d.Paragraphs(28)
$pStart = $para.Range.Start
$pEnd = $para.Range.End

# sanity check - only proceed on the paragraph that really holds the text
$pText = $para.Range.Text
if ($pText -like "*::Ask GPA only if the student is new::*") {

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>::</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Ask GPA only if the student is new::</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # target range = whole paragraph text minus the trailing paragraph mark
    $target = $d.Range($pStart, $pEnd - 1)
    $target.InsertXML($xmlFrag)
}

Write-Host "done"
